$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Datos actualizados" timestamp banner in A1
$ws.Range("A1").Value = "Datos actualizados a 19 de Mayo de 2020 a las 18:35"

# Row 4
$ws.Range("B4").Value = 1557770
$ws.Range("C4").Value = 7476
$ws.Range("D4").Value = 360058
$ws.Range("E4").Value = 1105234
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 497
$ws.Range("H4").Value = 92478

# Row 9
$ws.Range("B9").Value = 226699
$ws.Range("C9").Value = 813
$ws.Range("D9").Value = 129401
$ws.Range("E9").Value = 65129
$ws.Range("F9").Value = 0
$ws.Range("G9").Value = 162
$ws.Range("H9").Value = 32169

# Row 11
$ws.Range("B11").Value = 177620
$ws.Range("C11").Value = 331
$ws.Range("D11").Value = 155700
$ws.Range("E11").Value = 13770
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = 27
$ws.Range("H11").Value = 8150

# Row 12
$ws.Range("B12").Value = 151615
$ws.Range("C12").Value = 1022
$ws.Range("D12").Value = 112895
$ws.Range("E12").Value = 34521
$ws.Range("F12").Value = 0
$ws.Range("G12").Value = 28
$ws.Range("H12").Value = 4199

# Row 17
$ws.Range("B17").Value = 78500
$ws.Range("C17").Value = 428
$ws.Range("D17").Value = 39488
$ws.Range("E17").Value = 33154
$ws.Range("F17").Value = 0
$ws.Range("G17").Value = 16
$ws.Range("H17").Value = 5858

# Row 55
$ws.Range("B55").Value = 7377
$ws.Range("C55").Value = 176
$ws.Range("D55").Value = 3746
$ws.Range("E55").Value = 3070
$ws.Range("F55").Value = 0
$ws.Range("G55").Value = 6
$ws.Range("H55").Value = 561

# Row 58
$ws.Range("A58").Value = "Marruecos"
$ws.Range("B58").Value = 7023
$ws.Range("C58").Value = 71
$ws.Range("D58").Value = 3901
$ws.Range("E58").Value = 2929
$ws.Range("F58").Value = 0
$ws.Range("G58").Value = 1
$ws.Range("H58").Value = 193

# Row 59
$ws.Range("A59").Value = "Malasia"
$ws.Range("B59").Value = 6978
$ws.Range("C59").Value = 37
$ws.Range("D59").Value = 5646
$ws.Range("E59").Value = 1218
$ws.Range("F59").Value = 0
$ws.Range("G59").Value = 1
$ws.Range("H59").Value = 114

# Row 114
$ws.Range("A114").Value = "Paraguay"
$ws.Range("B114").Value = 829
$ws.Range("C114").Value = 41
$ws.Range("D114").Value = 230
$ws.Range("E114").Value = 588
$ws.Range("F114").Value = 0
$ws.Range("G114").Value = 0
$ws.Range("H114").Value = 11

# Row 115
$ws.Range("A115").Value = "Burkina Faso"
$ws.Range("B115").Value = 796
$ws.Range("C115").Value = 0
$ws.Range("D115").Value = 652
$ws.Range("E115").Value = 93
$ws.Range("F115").Value = 0
$ws.Range("G115").Value = 0
$ws.Range("H115").Value = 51

# Row 121
$ws.Range("B121").Value = 707
$ws.Range("C121").Value = 6
$ws.Range("D121").Value = 456
$ws.Range("E121").Value = 239
$ws.Range("F121").Value = 0
$ws.Range("G121").Value = 0
$ws.Range("H121").Value = 12

# Row 126
$ws.Range("A126").Value = "Republica del Chad"
$ws.Range("B126").Value = 545
$ws.Range("C126").Value = 26
$ws.Range("D126").Value = 139
$ws.Range("E126").Value = 350
$ws.Range("F126").Value = 0
$ws.Range("G126").Value = 3
$ws.Range("H126").Value = 56

# Row 127
$ws.Range("A127").Value = "Haiti"
$ws.Range("B127").Value = 533
$ws.Range("C127").Value = 0
$ws.Range("D127").Value = 21
$ws.Range("E127").Value = 491
$ws.Range("F127").Value = 0
$ws.Range("G127").Value = 0
$ws.Range("H127").Value = 21

# Row 128
$ws.Range("A128").Value = "Jamaica"
$ws.Range("B128").Value = 520
$ws.Range("C128").Value = 0
$ws.Range("D128").Value = 131
$ws.Range("E128").Value = 380
$ws.Range("F128").Value = 0
$ws.Range("G128").Value = 0
$ws.Range("H128").Value = 9

# Row 152
$ws.Range("B152").Value = 208
$ws.Range("C152").Value = 3
$ws.Range("D152").Value = 87
$ws.Range("E152").Value = 119
$ws.Range("F152").Value = 0
$ws.Range("G152").Value = 0
$ws.Range("H152").Value = 2

# Row 175
$ws.Range("B175").Value = 68
$ws.Range("C175").Value = 3
$ws.Range("D175").Value = 35
$ws.Range("E175").Value = 30
$ws.Range("F175").Value = 0
$ws.Range("G175").Value = 0
$ws.Range("H175").Value = 3
